$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-06-08 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-09 Sunday", 2)
$d.Content.Find.Execute("682÷3=227, 1", $true, $false, $false, $false, $false, $true, 1, $false, "897÷8=112, 1", 2)
$d.Content.Find.Execute("357÷9=39, 6", $true, $false, $false, $false, $false, $true, 1, $false, "610÷3=203, 1", 2)
$d.Content.Find.Execute("115÷5=23, 0", $true, $false, $false, $false, $false, $true, 1, $false, "337÷8=42, 1", 2)
$d.Content.Find.Execute("773÷6=128, 5", $true, $false, $false, $false, $false, $true, 1, $false, "703÷3=234, 1", 2)
$d.Content.Find.Execute("853÷8=106, 5", $true, $false, $false, $false, $false, $true, 1, $false, "822÷5=164, 2", 2)
$d.Content.Find.Execute("276÷2=138, 0", $true, $false, $false, $false, $false, $true, 1, $false, "512÷5=102, 2", 2)
$d.Content.Find.Execute("279÷6=46, 3", $true, $false, $false, $false, $false, $true, 1, $false, "651÷7=93, 0", 2)
$d.Content.Find.Execute("556÷4=139, 0", $true, $false, $false, $false, $false, $true, 1, $false, "195÷5=39, 0", 2)
$d.Content.Find.Execute("289÷8=36, 1", $true, $false, $false, $false, $false, $true, 1, $false, "711÷8=88, 7", 2)
$d.Content.Find.Execute("210÷8=26, 2", $true, $false, $false, $false, $false, $true, 1, $false, "994÷9=110, 4", 2)
$d.Content.Find.Execute("452÷3=150, 2", $true, $false, $false, $false, $false, $true, 1, $false, "825÷2=412, 1", 2)
$d.Content.Find.Execute("200÷4=50, 0", $true, $false, $false, $false, $false, $true, 1, $false, "848÷2=424, 0", 2)
$d.Content.Find.Execute("401÷4=100, 1", $true, $false, $false, $false, $false, $true, 1, $false, "163÷5=32, 3", 2)
$d.Content.Find.Execute("884÷9=98, 2", $true, $false, $false, $false, $false, $true, 1, $false, "121÷6=20, 1", 2)
$d.Content.Find.Execute("776÷5=155, 1", $true, $false, $false, $false, $false, $true, 1, $false, "558÷3=186, 0", 2)
$d.Content.Find.Execute("923÷6=153, 5", $true, $false, $false, $false, $false, $true, 1, $false, "287÷4=71, 3", 2)
$d.Content.Find.Execute("854÷2=427, 0", $true, $false, $false, $false, $false, $true, 1, $false, "226÷8=28, 2", 2)
$d.Content.Find.Execute("853÷6=142, 1", $true, $false, $false, $false, $false, $true, 1, $false, "377÷6=62, 5", 2)
$d.Content.Find.Execute("174÷3=58, 0", $true, $false, $false, $false, $false, $true, 1, $false, "519÷8=64, 7", 2)
$d.Content.Find.Execute("453÷5=90, 3", $true, $false, $false, $false, $false, $true, 1, $false, "391÷7=55, 6", 2)
$d.Content.Find.Execute("236÷6=39, 2", $true, $false, $false, $false, $false, $true, 1, $false, "637÷2=318, 1", 2)
$d.Content.Find.Execute("684÷4=171, 0", $true, $false, $false, $false, $false, $true, 1, $false, "844÷3=281, 1", 2)
$d.Content.Find.Execute("578÷8=72, 2", $true, $false, $false, $false, $false, $true, 1, $false, "680÷2=340, 0", 2)
$d.Content.Find.Execute("938÷9=104, 2", $true, $false, $false, $false, $false, $true, 1, $false, "238÷8=29, 6", 2)
$d.Content.Find.Execute("721÷2=360, 1", $true, $false, $false, $false, $false, $true, 1, $false, "366÷8=45, 6", 2)
